{"js": "const body = context.document.body;\n\n// Locate the \"hymns: ...\" front-matter paragraph by its text so the\n// script does not depend on a hard-coded paragraph index.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet hymnsPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"hymns:\") === 0) {\n    hymnsPara = paragraphs.items[i];\n    break;\n  }\n}\nif (!hymnsPara) {\n  throw new Error('Could not find the \"hymns:\" paragraph.');\n}\n\n// Insert a brand-new paragraph right after it containing \"speaker:\".\n// A placeholder marker character is appended first so the bookmark we\n// add below lands on an ordinary mid-run position instead of exactly\n// at the paragraph-end boundary (collapsed bookmarks placed exactly at\n// a paragraph mark land in the wrong spot).\nconst MARKER = \"\\u0001\";\nconst speakerPara = hymnsPara.insertParagraph(\"speaker:\" + MARKER, Word.InsertLocation.after);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark from the end of the \"category:\" paragraph\n// to the end of the new \"speaker:\" paragraph (right after \"speaker:\",\n// before the paragraph mark) - mirroring where Word leaves it after the\n// user's last edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst markerRange = speakerPara.search(MARKER, { matchCase: true }).getFirst();\nmarkerRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Remove the placeholder marker character now that the bookmark is\n// anchored in place.\nmarkerRange.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"hymns: ...\" front-matter paragraph by its text so the script\n# isn't dependent on a hard-coded paragraph index.\n$hymnsPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"hymns:*\") { $hymnsPara = $p; break }\n}\nif ($null -eq $hymnsPara) {\n    throw \"Could not find the 'hymns:' paragraph.\"\n}\n\n# Insert a brand new paragraph right after \"hymns: []\" and give it the\n# text \"speaker:\". A placeholder marker character is appended first so\n# the bookmark we add below lands on an ordinary mid-run position\n# instead of exactly at the paragraph-end boundary.\n$hymnsPara.Range.InsertParagraphAfter()\n$speakerPara = $hymnsPara.Next()\n$speakerPara.Range.Text = \"speaker:\" + [char]1\n\n# Move the \"_GoBack\" bookmark from the end of the \"category:\" paragraph\n# to the end of the new \"speaker:\" paragraph (right after \"speaker:\",\n# before the paragraph mark), mirroring where Word leaves it after the\n# user's last edit.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$markerRange = $speakerPara.Range\n$markerPos = $markerRange.End - 1\n$markerRange = $d.Range($markerPos - 1, $markerPos)\n$d.Bookmarks.Add(\"_GoBack\", $markerRange)\n\n# Remove the placeholder marker character now that the bookmark is\n# anchored in place.\n$markerRange.Text = \"\"\n"}
